# WebForm User Assignment execution
#
# The web-form "run" sheet carries one row per assigned test user; column F
# (header "PN_Value") holds that user's phone number. Re-running the
# assignment picked a fresh batch of phone numbers for rows 2-10, so update
# those nine cells in place.
#
# The column is plain digit strings stored as text (not numbers), so the
# range is switched to a Text number format before the values are typed in
# -- otherwise Excel would helpfully reinterpret "9840001395" as a numeric
# value. The style is reset back to Normal afterwards so the cells keep
# their original look-and-feel (General format) once the text is locked in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$phoneRange = $ws.Range("F2:F10")
$phoneRange.NumberFormat = "@"

$ws.Range("F2").Value  = "9840001395"
$ws.Range("F3").Value  = "9840094108"
$ws.Range("F4").Value  = "9840095536"
$ws.Range("F5").Value  = "9840058613"
$ws.Range("F6").Value  = "9840013354"
$ws.Range("F7").Value  = "9840026470"
$ws.Range("F8").Value  = "9840034579"
$ws.Range("F9").Value  = "9840086122"
$ws.Range("F10").Value = "9840091847"

$phoneRange.Style = "Normal"
